# Adds a "2022-Q4" sheet (new quarterly fund-holdings breakdown) right
# after "总计", and records its aggregate row (13 holdings, 1.67 亿元)
# as the new first data row on "总计" - shifting the older quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet immediately after "总计" (sheet 1).
#    All the other quarter sheets (2022-Q2 … 2020-Q4) shift right by one,
#    which is exactly the sheetId/tab-order change in the diff.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# Match the page margins used by every other sheet in the workbook
# (points: 0.75in/0.75in/1in/1in/0.5in/0.5in).
$q4Sheet.PageSetup.LeftMargin = 54
$q4Sheet.PageSetup.RightMargin = 54
$q4Sheet.PageSetup.TopMargin = 72
$q4Sheet.PageSetup.BottomMargin = 72
$q4Sheet.PageSetup.HeaderMargin = 36
$q4Sheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Fill in the "2022-Q4" fund-holdings table.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4Sheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# index, code, name, scale, stockPosition, positionPct, marketValue, rank
$rows = @(
    @(0,  "005299", "万家成长优选灵活配置混合A",     "12.13", "91.35", "3.60", "0.4367", 9),
    @(1,  "005300", "万家成长优选灵活配置混合C",     "9.48",  "91.35", "3.60", "0.3413", 9),
    @(2,  "010694", "万家内需增长一年持有期混合",     "9.46",  "94.46", "3.03", "0.2866", 10),
    @(3,  "010611", "万家战略发展产业混合A",         "5.75",  "92.07", "3.36", "0.1932", 9),
    @(4,  "010612", "万家战略发展产业混合C",         "4.14",  "92.07", "3.36", "0.1391", 9),
    @(5,  "006132", "万家智造优势混合A",             "4.10",  "93.92", "3.27", "0.1341", 8),
    @(6,  "159851", "华宝中证金融科技主题ETF",       "1.98",  "98.27", "2.87", "0.0568", 9),
    @(7,  "006133", "万家智造优势混合C",             "0.78",  "93.92", "3.27", "0.0255", 8),
    @(8,  "516100", "华夏中证金融科技主题ETF",       "0.60",  "97.54", "2.85", "0.0171", 9),
    @(9,  "000354", "长盛城镇化主题混合",             "0.33",  "85.83", "4.97", "0.0164", 6),
    @(10, "015112", "长盛精选行业轮动混合A",         "0.17",  "53.70", "4.98", "0.0085", 4),
    @(11, "516860", "博时中证金融科技主题ETF",       "0.27",  "98.47", "2.88", "0.0078", 9),
    @(12, "015113", "长盛精选行业轮动混合C",         "0.07",  "53.70", "4.98", "0.0035", 4)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $row = $i + 2

    # Column A: plain numeric index.
    $q4Sheet.Cells.Item($row, 1).Value = $r[0]

    # Columns B, D, E, F, G are stored as TEXT in this workbook (e.g. fund
    # codes keep leading zeros, percentages keep trailing zeros) - force
    # text storage via NumberFormat "@" so Excel doesn't silently coerce
    # them to numbers, then restore the plain "Normal" style so no stray
    # number format lingers on the cell (matches the un-styled text cells
    # elsewhere in the workbook).
    foreach ($col in 2, 4, 5, 6, 7) {
        $q4Sheet.Cells.Item($row, $col).NumberFormat = "@"
    }
    $q4Sheet.Cells.Item($row, 2).Value = $r[1]
    $q4Sheet.Cells.Item($row, 3).Value = $r[2]
    $q4Sheet.Cells.Item($row, 4).Value = $r[3]
    $q4Sheet.Cells.Item($row, 5).Value = $r[4]
    $q4Sheet.Cells.Item($row, 6).Value = $r[5]
    $q4Sheet.Cells.Item($row, 7).Value = $r[6]
    foreach ($col in 2, 4, 5, 6, 7) {
        $q4Sheet.Cells.Item($row, $col).Style = "Normal"
    }

    # Column H: plain numeric rank.
    $q4Sheet.Cells.Item($row, 8).Value = $r[7]
}

# Copy the header / row-label formatting (bold + border, centered) from an
# existing quarter sheet so the new sheet matches house style exactly.
$styleSource = $wb.Worksheets.Item("2022-Q2")
$q4Sheet2 = $wb.Worksheets.Item("2022-Q4")

$styleSource.Range("B1:H1").Copy()
$q4Sheet2.Range("B1:H1").PasteSpecial(-4122)

$styleSource.Range("A2:A5").Copy()
$q4Sheet2.Range("A2:A14").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Record the new quarter's totals as the first data row of "总计",
#    pushing the existing quarters down one row each.
# ---------------------------------------------------------------------
$totalData = @(
    @("2022-Q4", 13, 1.67),
    @("2022-Q2", 4, 0.08),
    @("2022-Q1", 3, 0.73),
    @("2021-Q4", 6, 0.84),
    @("2021-Q3", 4, 1.63),
    @("2021-Q2", 6, 4.09),
    @("2021-Q1", 6, 3.39),
    @("2020-Q4", 4, 3.59)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $row = $i + 2
    $d = $totalData[$i]
    $totalSheet.Cells.Item($row, 1).Value = $i
    $totalSheet.Cells.Item($row, 2).Value = $d[0]
    $totalSheet.Cells.Item($row, 3).Value = $d[1]
    $totalSheet.Cells.Item($row, 4).Value = $d[2]
}

# Row 9 (2020-Q4) is brand new - copy column A's index-cell style
# (bold, bordered, centered) from the row above it so it matches the rest
# of the index column.
$totalSheet.Cells.Item(8, 1).Copy()
$totalSheet.Cells.Item(9, 1).PasteSpecial(-4122)
$totalSheet.Cells.Item(9, 1).Value = 7

Write-Output "2022-Q4 sheet inserted and 总计 updated"
